# Apply the "adding diff and summing amps to scale amp inputs" edit to the
# minipowersupply workbook's "Sheet1" worksheet (channel truth-table).
#
# Top block (rows 2-6): permute the U-/U+/A-/A+ labels and turn row 5's
# entries into plain numbers (5 and 5) instead of text labels.
#
# Bottom table (rows 9-14, Table1): permute the 3-bit code labels assigned
# to each channel, and turn row 10's codes into plain numbers (100 and 100)
# instead of text labels, while the remaining codes stay as zero-padded text
# (the leading apostrophe keeps Excel from reinterpreting "010" etc. as a
# number and dropping the leading zero).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Top block: channel -> U/A pin labels -------------------------------
$ws.Range("B2").Value = "A-"
$ws.Range("C2").Value = "A+"

$ws.Range("C3").Value = "A-"

$ws.Range("B4").Value = "U-"
$ws.Range("C4").Value = "U-"

$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 5

$ws.Range("B6").Value = "U+"

# --- Table1 (rows 9-14): channel -> A0..2 / A3..5 select codes ---------
$ws.Range("B10").Value = 100
$ws.Range("C10").Value = 100

$ws.Range("B11").Value = "'010"
$ws.Range("C11").Value = "'010"

$ws.Range("B12").Value = "'001"
$ws.Range("C12").Value = "'000"

$ws.Range("B13").Value = "'000"
$ws.Range("C13").Value = "'001"

$ws.Range("B14").Value = "'011"
$ws.Range("C14").Value = "'011"

# --- Selection left on B15 after editing --------------------------------
$ws.Activate()
$ws.Range("B15").Select()
